# Add a new "Referencing" worksheet at the end of the workbook, give it some
# content and a workbook-level defined name ("Named_reference") that points
# at its A1 cell, then reference that name from a formula on the new sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the last existing sheet so it lands at
# the end of the tab strip (and becomes the active sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Referencing"

# A1 holds a label; it will be referenced by the defined name.
$ws.Range("A1").Value = "Named reference"

# Workbook-scoped defined name pointing at Referencing!$A$1
$wb.Names.Add("Named_reference", "=Referencing!`$A`$1")

# A2 uses the defined name in a formula.
$ws.Range("A2").Formula = "=Named_reference"

# Match the original author's selection/active cell on the new sheet.
[void]$ws.Range("A2").Select()

# Restore the page margins to match the rest of the workbook's sheets.
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36
